$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-100 down to 6-101.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data point.
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value = "Maule"
$ws.Cells.Item(5, 4).Value = 44630
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100103
$ws.Cells.Item(5, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value = 100103002
$ws.Cells.Item(5, 10).Value = "Ciruela"
$ws.Cells.Item(5, 11).Value = "Angeleno"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 8000
$ws.Cells.Item(5, 15).Value = 8000
$ws.Cells.Item(5, 16).Value = 8000
$ws.Cells.Item(5, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(5, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value = 444
$ws.Cells.Item(5, 20).Value = 18
